# Scheduled runner update: refresh computed market-price/profit figures
# (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6138.3076
$ws.Range("I33").Value = 7875.4
$ws.Range("K33").Value = 7875.4
$ws.Range("M33").Value = -7646.4
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -10968
$ws.Range("H86").Value = 4992.7896
$ws.Range("I86").Value = 4057
$ws.Range("J86").Value = 5327
$ws.Range("K86").Value = 4057
$ws.Range("L86").Value = 5327
$ws.Range("M86").Value = -2934
$ws.Range("N86").Value = -7573
$ws.Range("H89").Value = 4992.7896
$ws.Range("I89").Value = 4057
$ws.Range("J89").Value = 5327
$ws.Range("K89").Value = 20285
$ws.Range("L89").Value = 26635
$ws.Range("M89").Value = -14669
$ws.Range("N89").Value = -37867
$ws.Range("H106").Value = 17499.688
$ws.Range("I106").Value = 17856.928
$ws.Range("K106").Value = 17856.928
$ws.Range("M106").Value = -17225.928
$ws.Range("H112").Value = 6522.227
$ws.Range("I112").Value = 2245
$ws.Range("J112").Value = 7197.579
$ws.Range("K112").Value = 6735
$ws.Range("L112").Value = 21592.737
$ws.Range("M112").Value = -5627
$ws.Range("N112").Value = -23808.737
$ws.Range("H113").Value = 8478
$ws.Range("J113").Value = 8222.5
$ws.Range("L113").Value = 8222.5
$ws.Range("N113").Value = -14730.5
$ws.Range("H132").Value = 13336388
$ws.Range("I132").Value = 14928421
$ws.Range("K132").Value = 44785263
$ws.Range("M132").Value = -44782733
$ws.Range("H137").Value = 78473.83
$ws.Range("I137").Value = 137715.23
$ws.Range("J137").Value = 1460
$ws.Range("K137").Value = 413145.6900000001
$ws.Range("L137").Value = 4380
$ws.Range("M137").Value = -410595.6900000001
$ws.Range("N137").Value = -9480
$ws.Range("H138").Value = 4913.796
$ws.Range("I138").Value = 4552
$ws.Range("J138").Value = 4984.39
$ws.Range("K138").Value = 13656
$ws.Range("L138").Value = 14953.17
$ws.Range("M138").Value = -8516
$ws.Range("N138").Value = -25233.17
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9854.125
$ws.Range("I61").Value = 10178.066
$ws.Range("K61").Value = 10178.066
$ws.Range("M61").Value = -9966.066000000001
$ws.Range("H74").Value = 55225.31
$ws.Range("I74").Value = 10878.275
$ws.Range("J74").Value = 183831.7
$ws.Range("K74").Value = 10878.275
$ws.Range("L74").Value = 183831.7
$ws.Range("M74").Value = -10004.275
$ws.Range("N74").Value = -185579.7
$ws.Range("H77").Value = 55225.31
$ws.Range("I77").Value = 10878.275
$ws.Range("J77").Value = 183831.7
$ws.Range("K77").Value = 54391.375
$ws.Range("L77").Value = 919158.5
$ws.Range("M77").Value = -50023.375
$ws.Range("N77").Value = -927894.5
$ws.Range("H97").Value = 60733296
$ws.Range("I97").Value = 73740220
$ws.Range("K97").Value = 73740220
$ws.Range("M97").Value = -73739724
$ws.Range("H110").Value = 3473045.2
$ws.Range("I110").Value = 5556290
$ws.Range("K110").Value = 5556290
$ws.Range("M110").Value = -5554245
$ws.Range("H124").Value = 63475.668
$ws.Range("J124").Value = 63475.668
$ws.Range("L124").Value = 63475.668
$ws.Range("N124").Value = -73295.66800000001
$ws.Range("H132").Value = 5402.024
$ws.Range("I132").Value = 5649.385
$ws.Range("K132").Value = 16948.155
$ws.Range("M132").Value = -14418.155
$ws.Range("H136").Value = 9854.125
$ws.Range("I136").Value = 10178.066
$ws.Range("K136").Value = 30534.198
$ws.Range("M136").Value = -27984.198
$ws.Range("H139").Value = 322008.88
$ws.Range("J139").Value = 309248.78
$ws.Range("L139").Value = 309248.78
$ws.Range("N139").Value = -319528.78
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 8737
$ws.Range("I54").Value = 8737
$ws.Range("K54").Value = 8737
$ws.Range("M54").Value = -8253
$ws.Range("H122").Value = 110363
$ws.Range("J122").Value = 110363
$ws.Range("L122").Value = 110363
$ws.Range("N122").Value = -120163
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15404.419
$ws.Range("I31").Value = 5424.875
$ws.Range("K31").Value = 5424.875
$ws.Range("M31").Value = -5129.875
$ws.Range("H34").Value = 15404.419
$ws.Range("I34").Value = 5424.875
$ws.Range("K34").Value = 5424.875
$ws.Range("M34").Value = -5222.875
$ws.Range("H99").Value = 3359.6553
$ws.Range("I99").Value = 3191.3333
$ws.Range("J99").Value = 3635.0908
$ws.Range("K99").Value = 3191.3333
$ws.Range("L99").Value = 3635.0908
$ws.Range("M99").Value = -1693.3333
$ws.Range("N99").Value = -6631.0908
$ws.Range("H126").Value = 3359.6553
$ws.Range("I126").Value = 3191.3333
$ws.Range("J126").Value = 3635.0908
$ws.Range("K126").Value = 9573.999899999999
$ws.Range("L126").Value = 10905.2724
$ws.Range("M126").Value = -7103.999899999999
$ws.Range("N126").Value = -15845.2724
$ws.Range("H134").Value = 8341.314
$ws.Range("I134").Value = 5946.2085
$ws.Range("K134").Value = 17838.6255
$ws.Range("M134").Value = -15303.6255
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5574040.5
$ws.Range("I4").Value = 6784423.5
$ws.Range("K4").Value = 20353270.5
$ws.Range("M4").Value = -20353158.5
$ws.Range("H94").Value = 9714.286
$ws.Range("I94").Value = 8000
$ws.Range("K94").Value = 24000
$ws.Range("M94").Value = -23324
$ws.Range("H109").Value = 4685.75
$ws.Range("J109").Value = 4698
$ws.Range("L109").Value = 14094
$ws.Range("N109").Value = -16174
$ws.Range("H113").Value = 2773.9644
$ws.Range("I113").Value = 4878.5713
$ws.Range("K113").Value = 14635.7139
$ws.Range("M113").Value = -12465.7139
$ws.Range("H131").Value = 3780.6428
$ws.Range("I131").Value = 1105
$ws.Range("K131").Value = 3315
$ws.Range("M131").Value = 1725
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 605650.2
$ws.Range("J92").Value = 605650.2
$ws.Range("L92").Value = 605650.2
$ws.Range("N92").Value = -609394.2
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H94").Value = 22000
$ws.Range("J94").Value = 22000
$ws.Range("L94").Value = 22000
$ws.Range("N94").Value = -23352
$ws.Range("H97").Value = 2650947.5
$ws.Range("I97").Value = 4768106.5
$ws.Range("K97").Value = 4768106.5
$ws.Range("M97").Value = -4767610.5
$ws.Range("H126").Value = 8589758
$ws.Range("I126").Value = 9093787
$ws.Range("K126").Value = 27281361
$ws.Range("M126").Value = -27278891
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2938.6924
$ws.Range("I7").Value = 1887.3043
$ws.Range("K7").Value = 1887.3043
$ws.Range("M7").Value = -1775.3043
$ws.Range("H56").Value = 2283.6667
$ws.Range("I56").Value = 2283.6667
$ws.Range("K56").Value = 2283.6667
$ws.Range("M56").Value = -1592.6667
$ws.Range("H100").Value = 5073
$ws.Range("I100").Value = 3467.6667
$ws.Range("J100").Value = 5675
$ws.Range("K100").Value = 3467.6667
$ws.Range("L100").Value = 5675
$ws.Range("M100").Value = -2926.6667
$ws.Range("N100").Value = -6757
$ws.Range("H126").Value = 2938.6924
$ws.Range("I126").Value = 1887.3043
$ws.Range("K126").Value = 5661.9129
$ws.Range("M126").Value = -3191.9129
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3334.3333
$ws.Range("J17").Value = 3334.3333
$ws.Range("L17").Value = 3334.3333
$ws.Range("N17").Value = -3678.3333
$ws.Range("H23").Value = 5453
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -271
$ws.Range("H113").Value = 807.75
$ws.Range("I113").Value = 853.5
$ws.Range("J113").Value = 722.7857
$ws.Range("K113").Value = 2560.5
$ws.Range("L113").Value = 2168.3571
$ws.Range("M113").Value = -390.5
$ws.Range("N113").Value = -6508.3571
